# Updated cryptos list on Thu Nov 23 02:34:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# Row 2 - Bitcoin
Set-TextValue "D2" "37.388.09"
$ws.Range("E2").Value = "  +3.37%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.060.40"
$ws.Range("E3").Value = "  +4.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "237.27"
$ws.Range("E5").Value = "  +3.40%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +3.57%  "

# Row 7 - Solana
Set-TextValue "D7" "58.20"
$ws.Range("E7").Value = "  +8.63%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +5.12%  "

# Row 10 - OKB
Set-TextValue "D10" "57.89"
$ws.Range("E10").Value = "  +1.50%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.50%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +3.94%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.364.58"
$ws.Range("E13").Value = "  +4.80%  "

# Row 14 - Chainlink
Set-TextValue "D14" "14.38"
$ws.Range("E14").Value = "  +4.84%  "

# Row 15 - Avalanche
Set-TextValue "D15" "21.28"
$ws.Range("E15").Value = "  +8.12%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.779"
$ws.Range("E16").Value = "  +4.09%  "

# Row 17 - Polkadot
Set-TextValue "D17" "5.21"
$ws.Range("E17").Value = "  +4.60%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.070.78"
$ws.Range("E18").Value = "  +5.27%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "37.572.32"
$ws.Range("E19").Value = "  +4.12%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.16"
$ws.Range("E20").Value = "  +22.92%  "

# Row 21 - Litecoin
Set-TextValue "D21" "68.98"
$ws.Range("E21").Value = "  +2.71%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0813"
$ws.Range("E22").Value = "  +1.47%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "225.34"
$ws.Range("E23").Value = "  +2.57%  "

# Row 24 - Dai
Set-TextValue "D24" "1.00"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +6.42%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  +2.83%  "

# Row 27 - Monero
Set-TextValue "D27" "163.76"
$ws.Range("E27").Value = "  +2.75%  "

# Row 28 - Cosmos
Set-TextValue "D28" "8.88"
$ws.Range("E28").Value = "  +5.06%  "

# Row 29 - ImmutableX
$ws.Range("E29").Value = "  +11.12%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "19.22"
$ws.Range("E30").Value = "  +3.50%  "

# Row 31 - Kaspa
Set-TextValue "D31" "0.127"
$ws.Range("E31").Value = "  +5.19%  "

# Row 32 - Stellar
Set-TextValue "D32" "0.119"
$ws.Range("E32").Value = "  +2.78%  "

# Row 33 - Filecoin
Set-TextValue "D33" "4.50"
$ws.Range("E33").Value = "  +4.59%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0628"
$ws.Range("E34").Value = "  +5.13%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +12.92%  "

# Row 36 - InternetComputer(DFINITY)
Set-TextValue "D36" "4.48"
$ws.Range("E36").Value = "  +6.63%  "

# Row 37 - BinanceUSD
$ws.Range("E37").Value = "  +0.18%  "

# Row 38 - RenderToken
Set-TextValue "D38" "3.36"
$ws.Range("E38").Value = "  +5.13%  "

# Row 39 - WEMIXToken
$ws.Range("E39").Value = "  +0.85%  "

# Row 40 - THORChain
Set-TextValue "D40" "5.85"
$ws.Range("E40").Value = "  +13.99%  "

# Row 41 & 42 swap: FTXToken/Cronos -> Cronos/FTXToken (with updated values)
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D41" "0.0981"
$ws.Range("E41").Value = "  +12.09%  "

$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D42" "4.53"
$ws.Range("E42").Value = "  +27.01%  "

# Row 43 - HuobiToken
Set-TextValue "D43" "2.96"
$ws.Range("E43").Value = "  -2.11%  "

# Row 44 & 45 swap: Maker/Aave -> Aave/Maker (with updated values)
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "97.74"
$ws.Range("E44").Value = "  +12.48%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D45" "1.483.59"
$ws.Range("E45").Value = "  +4.52%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  +6.27%  "

# Row 47 - TrustWalletToken
Set-TextValue "D47" "1.15"
$ws.Range("E47").Value = "  +7.63%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "15.96"
$ws.Range("E48").Value = "  +9.57%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +4.62%  "

# Row 50 - FraxShare
Set-TextValue "D50" "7.22"
$ws.Range("E50").Value = "  +8.04%  "

# Row 51 - MXToken
Set-TextValue "D51" "2.94"
$ws.Range("E51").Value = "  +3.03%  "
